$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 348.72726
$ws.Range("I2").Value = 283.7
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 283.7
$ws.Range("L2").Value = 999
$ws.Range("M2").Value = -170.7
$ws.Range("N2").Value = -1225
$ws.Range("H9").Value = 6031.647
$ws.Range("I9").Value = 7822.077
$ws.Range("J9").Value = 212.75
$ws.Range("K9").Value = 7822.077
$ws.Range("L9").Value = 212.75
$ws.Range("M9").Value = -7653.077
$ws.Range("N9").Value = -550.75
$ws.Range("H29").Value = 2539.875
$ws.Range("I29").Value = 803.3333
$ws.Range("J29").Value = 7749.5
$ws.Range("K29").Value = 2409.9999
$ws.Range("L29").Value = 23248.5
$ws.Range("M29").Value = -2128.9999
$ws.Range("N29").Value = -23810.5
$ws.Range("H40").Value = 22225022
$ws.Range("I40").Value = 3500
$ws.Range("K40").Value = 3500
$ws.Range("M40").Value = -3325
$ws.Range("H86").Value = 11379.059
$ws.Range("I86").Value = 16421.857
$ws.Range("J86").Value = 7849.1
$ws.Range("K86").Value = 16421.857
$ws.Range("L86").Value = 7849.1
$ws.Range("M86").Value = -15298.857
$ws.Range("N86").Value = -10095.1
$ws.Range("H89").Value = 11379.059
$ws.Range("I89").Value = 16421.857
$ws.Range("J89").Value = 7849.1
$ws.Range("K89").Value = 82109.285
$ws.Range("L89").Value = 39245.5
$ws.Range("M89").Value = -76493.285
$ws.Range("N89").Value = -50477.5
$ws.Range("H100").Value = 1922.4
$ws.Range("I100").Value = 2063.3333
$ws.Range("K100").Value = 2063.3333
$ws.Range("M100").Value = -1522.3333
$ws.Range("H103").Value = 624.5
$ws.Range("J103").Value = 749
$ws.Range("L103").Value = 2247
$ws.Range("N103").Value = -3419
$ws.Range("H106").Value = 5613.8
$ws.Range("I106").Value = 5209.8335
$ws.Range("K106").Value = 5209.8335
$ws.Range("M106").Value = -4578.8335
$ws.Range("H125").Value = 7098440
$ws.Range("I125").Value = 1304297.2
$ws.Range("K125").Value = 11738674.8
$ws.Range("M125").Value = -11736214.8
$ws.Range("H132").Value = 2301.3333
$ws.Range("I132").Value = 2319.5454
$ws.Range("K132").Value = 6958.6362
$ws.Range("M132").Value = -4428.6362
$ws.Range("H135").Value = 32259146
$ws.Range("I135").Value = 1019.3077
$ws.Range("K135").Value = 9173.7693
$ws.Range("M135").Value = -6638.7693
$ws.Range("H138").Value = 6460.775
$ws.Range("J138").Value = 6446.237
$ws.Range("L138").Value = 19338.711
$ws.Range("N138").Value = -29618.711
$ws.Range("H139").Value = 112499.5
$ws.Range("J139").Value = 112499.5
$ws.Range("L139").Value = 112499.5
$ws.Range("N139").Value = -122779.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14037.678
$ws.Range("I32").Value = 7969.185
$ws.Range("J32").Value = 55000
$ws.Range("K32").Value = 7969.185
$ws.Range("L32").Value = 55000
$ws.Range("M32").Value = -7682.185
$ws.Range("N32").Value = -55574
$ws.Range("H36").Value = 17493.2
$ws.Range("I36").Value = 14364.5
$ws.Range("K36").Value = 14364.5
$ws.Range("M36").Value = -14018.5
$ws.Range("H61").Value = 33338400
$ws.Range("I61").Value = 35719356
$ws.Range("K61").Value = 35719356
$ws.Range("M61").Value = -35719144
$ws.Range("H97").Value = 1048.4762
$ws.Range("I97").Value = 964
$ws.Range("K97").Value = 964
$ws.Range("M97").Value = -468
$ws.Range("H110").Value = 1849.5264
$ws.Range("I110").Value = 1555.5
$ws.Range("J110").Value = 2063.3635
$ws.Range("K110").Value = 1555.5
$ws.Range("L110").Value = 2063.3635
$ws.Range("M110").Value = 489.5
$ws.Range("N110").Value = -6153.363499999999
$ws.Range("H122").Value = 2256
$ws.Range("I122").Value = 2256
$ws.Range("K122").Value = 6768
$ws.Range("M122").Value = -4318
$ws.Range("H136").Value = 33338400
$ws.Range("I136").Value = 35719356
$ws.Range("K136").Value = 107158068
$ws.Range("M136").Value = -107155518
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4147.5
$ws.Range("I86").Value = 4147.5
$ws.Range("K86").Value = 4147.5
$ws.Range("M86").Value = -3024.5
$ws.Range("H89").Value = 4147.5
$ws.Range("I89").Value = 4147.5
$ws.Range("K89").Value = 20737.5
$ws.Range("M89").Value = -15121.5
$ws.Range("H103").Value = 30417.375
$ws.Range("I103").Value = 10868
$ws.Range("K103").Value = 10868
$ws.Range("M103").Value = -9696
$ws.Range("H105").Value = 3573986.2
$ws.Range("I105").Value = 5953810.5
$ws.Range("K105").Value = 5953810.5
$ws.Range("M105").Value = -5952063.5
$ws.Range("H132").Value = 1835.8594
$ws.Range("I132").Value = 1710.5172
$ws.Range("J132").Value = 3047.5
$ws.Range("K132").Value = 5131.5516
$ws.Range("L132").Value = 9142.5
$ws.Range("M132").Value = -2601.5516
$ws.Range("N132").Value = -14202.5
$ws.Range("H134").Value = 2483.5518
$ws.Range("I134").Value = 2593.92
$ws.Range("K134").Value = 7781.76
$ws.Range("M134").Value = -5246.76
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 41284.44
$ws.Range("I5").Value = 50918.05
$ws.Range("J5").Value = 2750
$ws.Range("K5").Value = 152754.15
$ws.Range("L5").Value = 8250
$ws.Range("M5").Value = -152642.15
$ws.Range("N5").Value = -8474
$ws.Range("H12").Value = 385.88235
$ws.Range("I12").Value = 452.14285
$ws.Range("J12").Value = 339.5
$ws.Range("K12").Value = 1356.42855
$ws.Range("L12").Value = 1018.5
$ws.Range("M12").Value = -1183.42855
$ws.Range("N12").Value = -1364.5
$ws.Range("H46").Value = 257.2
$ws.Range("J46").Value = 368
$ws.Range("L46").Value = 1104
$ws.Range("N46").Value = -1286
$ws.Range("H55").Value = 700
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H68").Value = 734.2
$ws.Range("I68").Value = 703.625
$ws.Range("J68").Value = 856.5
$ws.Range("K68").Value = 2110.875
$ws.Range("L68").Value = 2569.5
$ws.Range("M68").Value = -1299.875
$ws.Range("N68").Value = -4191.5
$ws.Range("H71").Value = 734.2
$ws.Range("I71").Value = 703.625
$ws.Range("J71").Value = 856.5
$ws.Range("K71").Value = 6332.625
$ws.Range("L71").Value = 7708.5
$ws.Range("M71").Value = -2276.625
$ws.Range("N71").Value = -15820.5
$ws.Range("H135").Value = 41284.44
$ws.Range("I135").Value = 50918.05
$ws.Range("J135").Value = 2750
$ws.Range("K135").Value = 458262.45
$ws.Range("L135").Value = 24750
$ws.Range("M135").Value = -455727.45
$ws.Range("N135").Value = -29820
$ws.Range("H138").Value = 2970.6667
$ws.Range("J138").Value = 3900
$ws.Range("L138").Value = 11700
$ws.Range("N138").Value = -21980
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2364.1667
$ws.Range("J80").Value = 2257
$ws.Range("L80").Value = 2257
$ws.Range("N80").Value = -4253
$ws.Range("H83").Value = 2364.1667
$ws.Range("J83").Value = 2257
$ws.Range("L83").Value = 11285
$ws.Range("N83").Value = -21269
$ws.Range("H97").Value = 1132.125
$ws.Range("I97").Value = 920
$ws.Range("K97").Value = 920
$ws.Range("M97").Value = -424
$ws.Range("H122").Value = 175154.42
$ws.Range("I122").Value = 175154.42
$ws.Range("K122").Value = 525463.26
$ws.Range("M122").Value = -523013.26
$ws.Range("H126").Value = 7446.2104
$ws.Range("I126").Value = 6042.8667
$ws.Range("K126").Value = 18128.6001
$ws.Range("M126").Value = -15658.6001
$ws.Range("H132").Value = 7154.0625
$ws.Range("I132").Value = 6821.4
$ws.Range("K132").Value = 20464.2
$ws.Range("M132").Value = -17934.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2864
$ws.Range("H46").Value = 1713.4615
$ws.Range("I46").Value = 1598.9
$ws.Range("J46").Value = 2095.3333
$ws.Range("K46").Value = 1598.9
$ws.Range("L46").Value = 2095.3333
$ws.Range("M46").Value = -1410.9
$ws.Range("N46").Value = -2471.3333
$ws.Range("H82").Value = 1988.7646
$ws.Range("I82").Value = 1898.75
$ws.Range("K82").Value = 1898.75
$ws.Range("M82").Value = -1537.75
$ws.Range("H85").Value = 1988.7646
$ws.Range("I85").Value = 1898.75
$ws.Range("K85").Value = 1898.75
$ws.Range("M85").Value = -650.75
$ws.Range("H122").Value = 5454.9375
$ws.Range("I122").Value = 5061.727
$ws.Range("K122").Value = 15185.181
$ws.Range("M122").Value = -12735.181
$ws.Range("H132").Value = 17699670
$ws.Range("I132").Value = 20007628
$ws.Range("J132").Value = 5331.6665
$ws.Range("K132").Value = 60022884
$ws.Range("L132").Value = 15994.9995
$ws.Range("M132").Value = -60020354
$ws.Range("N132").Value = -21054.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 65178.8
$ws.Range("J95").Value = 65178.8
$ws.Range("L95").Value = 65178.8
$ws.Range("N95").Value = -70670.8
$ws.Range("H107").Value = 3399.2
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 3000
$ws.Range("M107").Value = -1080
$ws.Range("H132").Value = 5533.2095
$ws.Range("I132").Value = 3375.875
$ws.Range("K132").Value = 10127.625
$ws.Range("M132").Value = -7597.625
